# Generate Report for Archive
#
# 1. Replace every occurrence of the status text "Ready for handoff" with
#    "In Translation" on every worksheet (Overview, zh-cn, de-de).
# 2. Narrow the "zh-cn"/"de-de" columns on the Overview sheet and the
#    "Status" column on the per-locale sheets from ~17.22 characters down
#    to ~13.41 characters.
#    (Excel's ColumnWidth COM property is quantized to a 1/6-character
#    pixel grid on save, so 12.5 is the input that lands the stored width
#    closest to the intended 13.4101845877511.)

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used) {
        if ($cell.Text -eq "Ready for handoff") {
            $cell.Value = "In Translation"
        }
    }
}

$newWidth = 12.5

# Overview sheet: columns E (zh-cn) and F (de-de)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth

# zh-cn / de-de sheets: column C (Status)
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Columns.Item(3).ColumnWidth = $newWidth

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Columns.Item(3).ColumnWidth = $newWidth
